# Add the rest of the league schedule: 11 new match rows (46-56), inserted
# right after the existing match rows (which end at row 54 / match 45),
# pushing the "Total" summary block down from rows 55-60 to rows 66-71.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert 11 blank rows where the new matches go. Everything currently
#    at/after row 55 (the Format/Total summary block) shifts down to 66.
# ---------------------------------------------------------------------
$ws.Rows("55:65").Insert() | Out-Null

# ---------------------------------------------------------------------
# 2. Stamp the new rows with the same alternating-border formatting used
#    by the existing match rows (rows 53 & 54 form the two-row pattern
#    that repeats all the way up the table). Copying just the formats
#    and pasting across the whole 11-row block tiles the A/B/C/D../T
#    styles correctly (s=1/s=5 alternation on column A, s=20 on C, etc).
# ---------------------------------------------------------------------
$ws.Range("A53:T54").Copy() | Out-Null
$ws.Range("A55:T65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Fill in match number (A), format id (B), match name (C) and the
#    six VLOOKUP/RANK prize formulas (D,G,J,M,P,S) for each new row.
# ---------------------------------------------------------------------
$matches = @(
  @(46, "MI vs DC"),
  @(47, "RR vs CSK"),
  @(48, "RCB vs PBKS"),
  @(49, "KKR vs SRH"),
  @(50, "DC vs CSK"),
  @(51, "RR vs MI"),
  @(52, "RCB vs SRH"),
  @(53, "CSK vs PBKS"),
  @(54, "KKR vs RR"),
  @(55, "SRH vs MI"),
  @(56, "RCB vs DC")
)

$row = 55
foreach ($m in $matches) {
  $ws.Cells.Item($row, 1).Value = $m[0]          # A: match number
  $ws.Cells.Item($row, 2).Value = 2               # B: format id
  $ws.Cells.Item($row, 3).Value = $m[1]           # C: match name

  foreach ($pair in @(@(4,5), @(7,8), @(10,11), @(13,14), @(16,17), @(19,20))) {
    $fcol = $pair[0]   # formula column (D,G,J,M,P,S)
    $scol = $pair[1]   # score column   (E,H,K,N,Q,T)
    $fAddr = $ws.Cells.Item($row, $fcol).Address($false, $false)
    $sAddr = $ws.Cells.Item($row, $scol).Address($false, $false)
    $eAddr = $ws.Cells.Item($row, 5).Address($false, $false)
    $hAddr = $ws.Cells.Item($row, 8).Address($false, $false)
    $kAddr = $ws.Cells.Item($row, 11).Address($false, $false)
    $nAddr = $ws.Cells.Item($row, 14).Address($false, $false)
    $qAddr = $ws.Cells.Item($row, 17).Address($false, $false)
    $tAddr = $ws.Cells.Item($row, 20).Address($false, $false)
    $rankArgs = "(`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row)"
    $formula = "=IF(ISERROR(VLOOKUP(RANK($sAddr, $rankArgs, 0),  `$A`$2:`$C`$7, `$B$row+1, FALSE)),`"`",VLOOKUP(RANK($sAddr, $rankArgs, 0),  `$A`$2:`$C`$7, `$B$row+1, FALSE))"
    $ws.Cells.Item($row, $fcol).Formula = $formula
  }

  $row = $row + 1
}

# ---------------------------------------------------------------------
# 4. The conditional-formatting rules that used to sit on the "Total"
#    row (E57,H57,K57,N57,Q57,T57) need to move down to row 68 along
#    with that row.
# ---------------------------------------------------------------------
function Move-CondFormat($oldAddr, $newAddr) {
  $oldRng = $ws.Range($oldAddr)
  $fcs = $oldRng.FormatConditions
  for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $fc.ModifyAppliesToRange($ws.Range($newAddr))
  }
}

Move-CondFormat "E57" "E68"
Move-CondFormat "H57" "H68"
Move-CondFormat "K57" "K68"
Move-CondFormat "N57" "N68"
Move-CondFormat "Q57" "Q68"
Move-CondFormat "T57" "T68"

# ---------------------------------------------------------------------
# 5. Nudge the view so the frozen pane / selection land near the new
#    bottom of the table, same as the saved workbook.
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 42
$win.ScrollColumn = 1
$ws.Range("U68").Select() | Out-Null
